# KPI daily-template pagination: replace the week of 2023-09-18..23 entries
# with the next week of entries (2023-10-09..14), and move the active
# selection to F9 (per the author's "kpi pagination" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "task"
$ws.Range("C1").Value = "time"

$ws.Range("A2").Value = "2023-10-09"
$ws.Range("B2").Value = "Perjalanan ke Jakarta"
$ws.Range("C2").Value = "08:00"

$ws.Range("A3").Value = "2023-10-09"
$ws.Range("B3").Value = "Review achievement AMPM/KPI To Do List & Do & Done week 40"
$ws.Range("C3").Value = "09:00"

$ws.Range("A4").Value = "2023-10-09"
$ws.Range("B4").Value = "Create AMPM/To Do List Teddy"
$ws.Range("C4").Value = "11:00"

$ws.Range("A5").Value = "2023-10-09"
$ws.Range("B5").Value = "Stock opname Gudang MMI Prima Center"
$ws.Range("C5").Value = "13:00"

$ws.Range("A6").Value = "2023-10-09"
$ws.Range("B6").Value = "Penelusuran selisih Trenly Tambun"
$ws.Range("C6").Value = "16:00"

$ws.Range("A7").Value = "2023-10-09"
$ws.Range("B7").Value = "Review LHP week 40"
$ws.Range("C7").Value = "17:00"

$ws.Range("A8").Value = "2023-10-10"
$ws.Range("B8").Value = "Review To Do List"
$ws.Range("C8").Value = "09:00"

$ws.Range("A9").Value = "2023-10-10"
$ws.Range("B9").Value = "Stock opname Gudang MMI Prima Center"
$ws.Range("C9").Value = "09:30"

$ws.Range("A10").Value = "2023-10-10"
$ws.Range("B10").Value = "Review LHP week 40"
$ws.Range("C10").Value = "16:00"

$ws.Range("A11").Value = "2023-10-11"
$ws.Range("B11").Value = "Morning briefing Retail"
$ws.Range("C11").Value = "08:30"

$ws.Range("A12").Value = "2023-10-11"
$ws.Range("B12").Value = "Review To Do List "
$ws.Range("C12").Value = "09:30"

$ws.Range("A13").Value = "2023-10-11"
$ws.Range("B13").Value = "Review audit CV Maju Technology dengan Joe di Ocean Space"
$ws.Range("C13").Value = "10:00"

$ws.Range("A14").Value = "2023-10-11"
$ws.Range("B14").Value = "Review audit CCTV & audit konsumen dengan Eni di Ocean Space"
$ws.Range("C14").Value = "15:00"

$ws.Range("A15").Value = "2023-10-11"
$ws.Range("B15").Value = "Review LHP week 40"
$ws.Range("C15").Value = "17:00"

$ws.Range("A16").Value = "2023-10-12"
$ws.Range("B16").Value = "Review To Do List"
$ws.Range("C16").Value = "08:30"

$ws.Range("A17").Value = "2023-10-12"
$ws.Range("B17").Value = "Review finding status dengan Faizal & Lula"
$ws.Range("C17").Value = "09:00"

$ws.Range("A18").Value = "2023-10-12"
$ws.Range("B18").Value = "appraisal dengan Fajar di CS"
$ws.Range("C18").Value = "13:00"

$ws.Range("A19").Value = "2023-10-12"
$ws.Range("B19").Value = "Review Telemarketing"
$ws.Range("C19").Value = "14:00"

$ws.Range("A20").Value = "2023-10-12"
$ws.Range("B20").Value = "Review audit Dolphin dengan Ridwan"
$ws.Range("C20").Value = "15:00"

$ws.Range("A21").Value = "2023-10-12"
$ws.Range("B21").Value = "LHP week 40"
$ws.Range("C21").Value = "16:00"

$ws.Range("A22").Value = "2023-10-13"
$ws.Range("B22").Value = "Review To Do List "
$ws.Range("C22").Value = "08:30"

$ws.Range("A23").Value = "2023-10-13"
$ws.Range("B23").Value = "LPJ Coordinator Internal Audit "
$ws.Range("C23").Value = "09:00"

$ws.Range("A24").Value = "2023-10-13"
$ws.Range("B24").Value = "Meeting CV SMJ "
$ws.Range("C24").Value = "14:00"

$ws.Range("A25").Value = "2023-10-13"
$ws.Range("B25").Value = "Meeting Kospin SMS"
$ws.Range("C25").Value = "16:00"

$ws.Range("A26").Value = "2023-10-14"
$ws.Range("B26").Value = "Review To Do List"
$ws.Range("C26").Value = "08:30"

$ws.Range("A27").Value = "2023-10-14"
$ws.Range("B27").Value = "Upgrade skill audit CV CS"
$ws.Range("C27").Value = "09:00"

$ws.Range("A28").Value = "2023-10-14"
$ws.Range("B28").Value = "Review LHP week 40"
$ws.Range("C28").Value = "14:00 "

$ws.Range("F9").Select()
